{"js": "const pairs = [\n  [\"233\u00d72=\", \"558\u00d73=\"],\n  [\"436\u00d73=\", \"788\u00d77=\"],\n  [\"921\u00d79=\", \"733\u00d75=\"],\n  [\"207\u00d74=\", \"152\u00d74=\"],\n  [\"988\u00d77=\", \"377\u00d76=\"],\n  [\"150\u00d75=\", \"973\u00d76=\"],\n  [\"638\u00d76=\", \"201\u00d77=\"],\n  [\"126\u00d78=\", \"483\u00d79=\"],\n  [\"858\u00d75=\", \"512\u00d79=\"],\n  [\"462\u00d79=\", \"479\u00d77=\"],\n  [\"341\u00d74=\", \"101\u00d73=\"],\n  [\"103\u00d76=\", \"157\u00d75=\"],\n  [\"868\u00d75=\", \"671\u00d76=\"],\n  [\"301\u00d74=\", \"612\u00d74=\"],\n  [\"233\u00d78=\", \"799\u00d73=\"],\n  [\"293\u00d79=\", \"231\u00d78=\"],\n  [\"758\u00d75=\", \"540\u00d72=\"],\n  [\"244\u00d76=\", \"892\u00d77=\"],\n  [\"179\u00d75=\", \"720\u00d72=\"],\n  [\"152\u00d76=\", \"867\u00d73=\"],\n  [\"616\u00d79=\", \"177\u00d78=\"],\n  [\"720\u00d77=\", \"233\u00d75=\"],\n  [\"621\u00d77=\", \"293\u00d72=\"],\n  [\"817\u00d72=\", \"903\u00d73=\"],\n  [\"198\u00d74=\", \"766\u00d76=\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  // Always sync after each replacement batch so later searches see the\n  // updated text (and so a later \"new\" value that happens to equal an\n  // earlier \"old\" value can't be matched twice).\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @('233\u00d72=', '558\u00d73='),\n    @('436\u00d73=', '788\u00d77='),\n    @('921\u00d79=', '733\u00d75='),\n    @('207\u00d74=', '152\u00d74='),\n    @('988\u00d77=', '377\u00d76='),\n    @('150\u00d75=', '973\u00d76='),\n    @('638\u00d76=', '201\u00d77='),\n    @('126\u00d78=', '483\u00d79='),\n    @('858\u00d75=', '512\u00d79='),\n    @('462\u00d79=', '479\u00d77='),\n    @('341\u00d74=', '101\u00d73='),\n    @('103\u00d76=', '157\u00d75='),\n    @('868\u00d75=', '671\u00d76='),\n    @('301\u00d74=', '612\u00d74='),\n    @('233\u00d78=', '799\u00d73='),\n    @('293\u00d79=', '231\u00d78='),\n    @('758\u00d75=', '540\u00d72='),\n    @('244\u00d76=', '892\u00d77='),\n    @('179\u00d75=', '720\u00d72='),\n    @('152\u00d76=', '867\u00d73='),\n    @('616\u00d79=', '177\u00d78='),\n    @('720\u00d77=', '233\u00d75='),\n    @('621\u00d77=', '293\u00d72='),\n    @('817\u00d72=', '903\u00d73='),\n    @('198\u00d74=', '766\u00d76='),\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.MatchSoundsLike = $false\n    $find.MatchAllWordForms = $false\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n\nWrite-Output \"replacements complete\"\n"}
